$wb = $excel.ActiveWorkbook

# Rename the three sheets (PROS-9738 - CCRU - New POS 2019 KPIs)
$wb.Worksheets.Item("PoS 2019 - Cinema - CAP ").Name  = "PoS 2019 - IC Cinema - CAP"
$wb.Worksheets.Item("PoS 2019 - Cinema - REG ").Name  = "PoS 2019 - IC Cinema - REG"
$wb.Worksheets.Item("PoS 2019 - FastFood ").Name      = "PoS 2019 - IC FastFood"

# Move the active tab/selected sheet from "IC Petroleum - REG" (index 12)
# to "IC HoReCa RestCafeTea" (index 6)
$wb.Worksheets.Item("PoS 2019 - IC HoReCa RestCafeTea").Activate()
